# Auto-generated Excel COM-interop script
# Updates numeric cell values in the Cactuar Profits workbook tables
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR per the authored diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 1500
$ws.Range("J10").Value = 1500
$ws.Range("L10").Value = 1500
$ws.Range("N10").Value = -2086
$ws.Range("H17").Value = 693.9474
$ws.Range("J17").Value = 693.9474
$ws.Range("L17").Value = 2081.8422
$ws.Range("N17").Value = -2417.8422
$ws.Range("H28").Value = 846.46155
$ws.Range("I28").Value = 368.5
$ws.Range("K28").Value = 368.5
$ws.Range("M28").Value = 116.5
$ws.Range("H86").Value = 1245814.2
$ws.Range("I86").Value = 1797822.2
$ws.Range("J86").Value = 3796.125
$ws.Range("K86").Value = 1797822.2
$ws.Range("L86").Value = 3796.125
$ws.Range("M86").Value = -1796699.2
$ws.Range("N86").Value = -6042.125
$ws.Range("H89").Value = 1245814.2
$ws.Range("I89").Value = 1797822.2
$ws.Range("J89").Value = 3796.125
$ws.Range("K89").Value = 8989111
$ws.Range("L89").Value = 18980.625
$ws.Range("M89").Value = -8983495
$ws.Range("N89").Value = -30212.625
$ws.Range("H100").Value = 1423
$ws.Range("I100").Value = 1423
$ws.Range("K100").Value = 1423
$ws.Range("M100").Value = -882
$ws.Range("H107").Value = 540.63635
$ws.Range("I107").Value = 494.1
$ws.Range("K107").Value = 494.1
$ws.Range("M107").Value = 1425.9
$ws.Range("H132").Value = 10004.258
$ws.Range("I132").Value = 3047.76
$ws.Range("K132").Value = 9143.280000000001
$ws.Range("M132").Value = -6613.280000000001
$ws.Range("H138").Value = 5947.1885
$ws.Range("I138").Value = 973.7857
$ws.Range("J138").Value = 7732.5127
$ws.Range("K138").Value = 2921.3571
$ws.Range("L138").Value = 23197.5381
$ws.Range("M138").Value = 2218.6429
$ws.Range("N138").Value = -33477.53810000001
$ws.Range("H141").Value = 6608.05
$ws.Range("I141").Value = 6166.421
$ws.Range("J141").Value = 14999
$ws.Range("K141").Value = 18499.263
$ws.Range("L141").Value = 44997
$ws.Range("M141").Value = -13319.263
$ws.Range("N141").Value = -55357

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3109.8235
$ws.Range("I45").Value = 2705
$ws.Range("K45").Value = 2705
$ws.Range("M45").Value = -2328
$ws.Range("H63").Value = 4250
$ws.Range("I63").Value = 4250
$ws.Range("K63").Value = 4250
$ws.Range("M63").Value = -3564
$ws.Range("H66").Value = 4250
$ws.Range("I66").Value = 4250
$ws.Range("K66").Value = 21250
$ws.Range("M66").Value = -17818
$ws.Range("H122").Value = 5859.231
$ws.Range("J122").Value = 9620.666999999999
$ws.Range("L122").Value = 28862.001
$ws.Range("N122").Value = -33762.001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2780.2222
$ws.Range("I20").Value = 2284.55
$ws.Range("K20").Value = 2284.55
$ws.Range("M20").Value = -2037.55
$ws.Range("H22").Value = 619.6
$ws.Range("J22").Value = 899.5
$ws.Range("L22").Value = 899.5
$ws.Range("N22").Value = -1245.5
$ws.Range("H134").Value = 6360.5
$ws.Range("I134").Value = 3934.5
$ws.Range("K134").Value = 11803.5
$ws.Range("M134").Value = -9268.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 4882.8335
$ws.Range("I10").Value = 4765.6665
$ws.Range("K10").Value = 4765.6665
$ws.Range("M10").Value = -4626.6665
$ws.Range("H13").Value = 450
$ws.Range("I13").Value = 450
$ws.Range("K13").Value = 450
$ws.Range("M13").Value = -311
$ws.Range("H31").Value = 4133
$ws.Range("J31").Value = 10000
$ws.Range("L31").Value = 10000
$ws.Range("N31").Value = -10590
$ws.Range("H34").Value = 4133
$ws.Range("J34").Value = 10000
$ws.Range("L34").Value = 10000
$ws.Range("N34").Value = -10404
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("H99").Value = 7785.357
$ws.Range("I99").Value = 3799.2
$ws.Range("K99").Value = 3799.2
$ws.Range("M99").Value = -2301.2
$ws.Range("H108").Value = 47300
$ws.Range("J108").Value = 47300
$ws.Range("L108").Value = 47300
$ws.Range("N108").Value = -54980
$ws.Range("H122").Value = 3113.037
$ws.Range("I122").Value = 1955.7222
$ws.Range("K122").Value = 5867.1666
$ws.Range("M122").Value = -3417.1666
$ws.Range("H126").Value = 7785.357
$ws.Range("I126").Value = 3799.2
$ws.Range("K126").Value = 11397.6
$ws.Range("M126").Value = -8927.599999999999
$ws.Range("H132").Value = 19629802
$ws.Range("I132").Value = 25659440
$ws.Range("K132").Value = 76978320
$ws.Range("M132").Value = -76975790
$ws.Range("H140").Value = 49904.707
$ws.Range("J140").Value = 49904.707
$ws.Range("L140").Value = 49904.707
$ws.Range("N140").Value = -60264.707
$ws.Range("N53").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 885.0769
$ws.Range("J5").Value = 949.1
$ws.Range("L5").Value = 2847.3
$ws.Range("N5").Value = -3071.3
$ws.Range("H132").Value = 1835.6
$ws.Range("J132").Value = 2499.5
$ws.Range("L132").Value = 22495.5
$ws.Range("N132").Value = -27555.5
$ws.Range("H135").Value = 885.0769
$ws.Range("J135").Value = 949.1
$ws.Range("L135").Value = 8541.9
$ws.Range("N135").Value = -13611.9

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3669005.5
$ws.Range("I70").Value = 9529314
$ws.Range("J70").Value = 6312.5
$ws.Range("K70").Value = 9529314
$ws.Range("L70").Value = 6312.5
$ws.Range("M70").Value = -9529044
$ws.Range("N70").Value = -6852.5
$ws.Range("H73").Value = 3669005.5
$ws.Range("I73").Value = 9529314
$ws.Range("J73").Value = 6312.5
$ws.Range("K73").Value = 9529314
$ws.Range("L73").Value = 6312.5
$ws.Range("M73").Value = -9528378
$ws.Range("N73").Value = -8184.5
$ws.Range("H102").Value = 14172472
$ws.Range("I102").Value = 18894360
$ws.Range("K102").Value = 18894360
$ws.Range("M102").Value = -18892738
$ws.Range("H113").Value = 2303.3333
$ws.Range("I113").Value = 2303.3333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2303.3333
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -133.3332999999998
$ws.Range("H122").Value = 4833.4165
$ws.Range("I122").Value = 1755.1666
$ws.Range("K122").Value = 5265.4998
$ws.Range("M122").Value = -2815.4998
$ws.Range("H126").Value = 3618.6333
$ws.Range("I126").Value = 2917.2104
$ws.Range("J126").Value = 4830.1816
$ws.Range("K126").Value = 8751.6312
$ws.Range("L126").Value = 14490.5448
$ws.Range("M126").Value = -6281.6312
$ws.Range("N126").Value = -19430.5448
$ws.Range("N113").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H122").Value = 57147340
$ws.Range("I122").Value = 100003630
$ws.Range("J122").Value = 14291049
$ws.Range("K122").Value = 300010890
$ws.Range("L122").Value = 42873147
$ws.Range("M122").Value = -300008440
$ws.Range("N122").Value = -42878047
$ws.Range("H136").Value = 3009.1428
$ws.Range("I136").Value = 2325.8948
$ws.Range("J136").Value = 9500
$ws.Range("K136").Value = 6977.6844
$ws.Range("L136").Value = 28500
$ws.Range("M136").Value = -4427.6844
$ws.Range("N136").Value = -33600
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 237.83333
$ws.Range("I4").Value = 77.90000000000001
$ws.Range("J4").Value = 1037.5
$ws.Range("K4").Value = 77.90000000000001
$ws.Range("L4").Value = 1037.5
$ws.Range("M4").Value = 35.09999999999999
$ws.Range("N4").Value = -1263.5
$ws.Range("H122").Value = 4752.091
$ws.Range("I122").Value = 4030.4443
$ws.Range("J122").Value = 7999.5
$ws.Range("K122").Value = 12091.3329
$ws.Range("L122").Value = 23998.5
$ws.Range("M122").Value = -9641.332900000001
$ws.Range("N122").Value = -28898.5
$ws.Range("H126").Value = 2129.3333
$ws.Range("I126").Value = 1194.5
$ws.Range("K126").Value = 3583.5
$ws.Range("M126").Value = -1113.5
$ws.Range("H132").Value = 42738024
$ws.Range("I132").Value = 7939078
$ws.Range("J132").Value = 83336790
$ws.Range("K132").Value = 23817234
$ws.Range("L132").Value = 250010370
$ws.Range("M132").Value = -23814704
$ws.Range("N132").Value = -250015430
